# Refresh the "cryptos" price/volume table with the latest scraped values.
# Mirrors the daily GitHub Actions update of cryptos.xlsx: prices and 1h volume
# deltas are refreshed per-coin; rows 22/23 (Uniswap/Polygon) also swap ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '66.810.45'
$ws.Range('E2').Value = '  -0.58%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '3.111.59'
$ws.Range('E3').Value = '  +0.05%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.02%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.84%  '

# Row 6: Solana
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.20'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.64%  '

# Row 7: USDC
$ws.Range('E7').Value = '  +0.05%  '

# Row 8: LidoStakedEther
$ws.Range('D8').Value = '3.108.47'
$ws.Range('E8').Value = '  +0.13%  '

# Row 9: XRP
$ws.Range('E9').Value = '  -0.65%  '

# Row 10: Toncoin
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.46'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.42%  '

# Row 11: Dogecoin
$ws.Range('E11').Value = '  -1.98%  '

# Row 12: Cardano
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.482'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.28%  '

# Row 13: ShibaInu
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000245'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.24%  '

# Row 14: Avalanche
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.27'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.59%  '

# Row 15: TRON
$ws.Range('E15').Value = '  -1.48%  '

# Row 16: WrappedliquidstakedEther2.0
$ws.Range('D16').Value = '3.628.45'
$ws.Range('E16').Value = '  +0.12%  '

# Row 17: WrappedBTC
$ws.Range('D17').Value = '66.786.67'
$ws.Range('E17').Value = '  -0.67%  '

# Row 18: Polkadot
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.14'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.46%  '

# Row 19: WrappedEther
$ws.Range('D19').Value = '3.108.91'
$ws.Range('E19').Value = '  -0.02%  '

# Row 20: Chainlink
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.97%  '

# Row 21: BitcoinCash
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '477.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.90%  '

# Row 22: Uniswap
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.714'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.63%  '

# Row 23: Polygon
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.99'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.71%  '

# Row 24: InternetComputer(DFINITY)
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.14%  '

# Row 25: Litecoin
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.99'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.16%  '

# Row 26: Fetch.AI
$ws.Range('E26').Value = '  -1.59%  '

# Row 27: RenderToken
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.03%  '

# Row 28: Dai
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.05%  '

# Row 29: ImmutableX
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.41'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.84%  '

# Row 30: NEARProtocol
$ws.Range('E30').Value = '  -2.36%  '

# Row 31: PancakeSwap
$ws.Range('E31').Value = '  -0.91%  '

# Row 32: EthereumClassic
$ws.Range('E32').Value = '  +0.07%  '

# Row 33: Hedera
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.115'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.32%  '

# Row 34: PEPE
$ws.Range('D34').Value = '0.0₃0938'
$ws.Range('E34').Value = '  -7.96%  '

# Row 35: FirstDigitalUSD
$ws.Range('E35').Value = '  -0.03%  '

# Row 36: Filecoin
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.86'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.19%  '

# Row 37: Mantle
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.975'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.94%  '

# Row 38: Arweave
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '47.17'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.09%  '

# Row 39: Stacks
$ws.Range('E39').Value = '  -0.62%  '

# Row 40: OKB
$ws.Range('E40').Value = '  -0.92%  '

# Row 41: TheGraph
$ws.Range('E41').Value = '  -2.33%  '

# Row 42: Kaspa
$ws.Range('E42').Value = '  -1.89%  '

# Row 43: Cosmos
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.71'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.22%  '

# Row 44: Maker
$ws.Range('D44').Value = '2.808.57'
$ws.Range('E44').Value = '  +1.18%  '

# Row 45: VeChain
$ws.Range('E45').Value = '  -2.65%  '

# Row 46: Bittensor
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '379.87'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.32%  '

# Row 47: dogwifhat
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.57'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -11.17%  '

# Row 48: Monero
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '136.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.73%  '

# Row 49: USDe
$ws.Range('E49').Value = '  +0.08%  '

# Row 50: InjectiveProtocol
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.73%  '

# Row 51: ThetaToken
$ws.Range('E51').Value = '  -2.29%  '
